$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $newValue
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '28.203.15'
Set-TextValue 'E2' '  -1.51%  '
Set-TextValue 'D3' '1.805.40'
Set-TextValue 'E3' '  +0.36%  '
Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  -0.07%  '
Set-TextValue 'D5' '316.88'
Set-TextValue 'E5' '  +1.08%  '
Set-TextValue 'E6' '  -0.06%  '
Set-TextValue 'D7' '0.5384'
Set-TextValue 'E7' '  +1.29%  '
Set-TextValue 'D8' '0.3786'
Set-TextValue 'E8' '  +0.53%  '
Set-TextValue 'D9' '0.07489'
Set-TextValue 'E9' '  -0.58%  '
Set-TextValue 'D10' '41.99'
Set-TextValue 'E10' '  -1.32%  '
Set-TextValue 'E11' '  -2.09%  '
Set-TextValue 'D12' '0.9999'
Set-TextValue 'E12' '  -0.08%  '
Set-TextValue 'D13' '6.213'
Set-TextValue 'E13' '  +0.08%  '
Set-TextValue 'D14' '20.58'
Set-TextValue 'E14' '  -3.04%  '
Set-TextValue 'D15' '7.396'
Set-TextValue 'E15' '  -1.04%  '
Set-TextValue 'D16' '1.803.55'
Set-TextValue 'E16' '  +0.46%  '
Set-TextValue 'D17' '90.01'
Set-TextValue 'E17' '  -0.64%  '
Set-TextValue 'D18' '0.00001067'
Set-TextValue 'E18' '  -0.30%  '
Set-TextValue 'D19' '0.06507'
Set-TextValue 'E19' '  +0.77%  '
Set-TextValue 'D20' '17.45'
Set-TextValue 'E20' '  +0.74%  '
Set-TextValue 'D21' '0.9997'
Set-TextValue 'E21' '  -0.04%  '
Set-TextValue 'D22' '5.939'
Set-TextValue 'E22' '  +0.13%  '
Set-TextValue 'D23' '28.225.43'
Set-TextValue 'D24' '11.21'
Set-TextValue 'E24' '  +0.19%  '
Set-TextValue 'D25' '2.089'
Set-TextValue 'E25' '  -0.36%  '
Set-TextValue 'D26' '156.06'
Set-TextValue 'E26' '  -3.07%  '
Set-TextValue 'D27' '20.55'
Set-TextValue 'E27' '  -0.03%  '
Set-TextValue 'D28' '2.010.96'
Set-TextValue 'E28' '  +0.54%  '
Set-TextValue 'D29' '2.337'
Set-TextValue 'E29' '  -2.25%  '
Set-TextValue 'D30' '122.23'
Set-TextValue 'E30' '  -1.18%  '
Set-TextValue 'B31' 'Stellar'
Set-TextValue 'C31' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D31' '0.1116'
Set-TextValue 'E31' '  +9.02%  '
Set-TextValue 'B32' 'ImmutableX'
Set-TextValue 'C32' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D32' '1.128'
Set-TextValue 'E32' '  -0.30%  '
Set-TextValue 'D33' '5.615'
Set-TextValue 'E33' '  -1.94%  '
Set-TextValue 'E34' '  +0.14%  '
Set-TextValue 'D35' '0.06975'
Set-TextValue 'E35' '  +6.49%  '
Set-TextValue 'D36' '0.2233'
Set-TextValue 'E36' '  -3.26%  '
Set-TextValue 'D37' '0.02305'
Set-TextValue 'E37' '  -0.84%  '
Set-TextValue 'D38' '5.103'
Set-TextValue 'E38' '  +0.55%  '
Set-TextValue 'D39' '8.475'
Set-TextValue 'E39' '  -4.04%  '
Set-TextValue 'E40' '  -2.78%  '
Set-TextValue 'D41' '0.6190'
Set-TextValue 'E41' '  -2.09%  '
Set-TextValue 'B42' 'TrustWalletToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D42' '1.179'
Set-TextValue 'E42' '  -2.14%  '
Set-TextValue 'B43' 'WEMIXTOKEN'
Set-TextValue 'C43' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D43' '1.429'
Set-TextValue 'E43' '  +2.53%  '
Set-TextValue 'D44' '13.43'
Set-TextValue 'E44' '  -0.55%  '
Set-TextValue 'D45' '3.688'
Set-TextValue 'E45' '  +0.56%  '
Set-TextValue 'D46' '0.5784'
Set-TextValue 'E46' '  -2.51%  '
Set-TextValue 'D47' '125.45'
Set-TextValue 'E47' '  -0.40%  '
Set-TextValue 'E48' '  +1.51%  '
Set-TextValue 'E49' '  -2.43%  '
Set-TextValue 'D50' '0.06828'
Set-TextValue 'E50' '  -1.54%  '
Set-TextValue 'D51' '72.07'
Set-TextValue 'E51' '  -1.62%  '
